$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D; old D:K data shifts right to F:M
$ws.Range("D:E").Insert()

# Copy number formatting from column F (the shifted old column D) into new D:E columns,
# restricted to the row blocks that actually had data in the original D:K columns
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = Dec-2018, E = Sep-2018)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 229400
$ws.Range("E8").Value = 218400
$ws.Range("D9").Value = 171600
$ws.Range("E9").Value = 165900
$ws.Range("D10").Value = 57800
$ws.Range("E10").Value = 52500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 77700
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 18200
$ws.Range("E15").Value = 15500
$ws.Range("D17").Value = 290300
$ws.Range("E17").Value = 202100
$ws.Range("D18").Value = -60900
$ws.Range("E18").Value = 16300
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = -42700
$ws.Range("E21").Value = 31900
$ws.Range("D22").Value = 16000
$ws.Range("E22").Value = 1600
$ws.Range("D23").Value = -76800
$ws.Range("E23").Value = 14800
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 1100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -77300
$ws.Range("E26").Value = 13700
$ws.Range("D27").Value = -77300
$ws.Range("E27").Value = 13700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -77300
$ws.Range("E33").Value = 13700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -77300
$ws.Range("E35").Value = 13700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 63600
$ws.Range("E41").Value = 86500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 162400
$ws.Range("E43").Value = 173100
$ws.Range("D44").Value = 91400
$ws.Range("E44").Value = 29600
$ws.Range("D45").Value = 15700
$ws.Range("E45").Value = 7000
$ws.Range("D46").Value = 333200
$ws.Range("E46").Value = 296200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 211600
$ws.Range("E48").Value = 257400
$ws.Range("D49").Value = 590000
$ws.Range("E49").Value = 151600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 6400
$ws.Range("E52").Value = 1100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1141200
$ws.Range("E54").Value = 706500
$ws.Range("D57").Value = 46100
$ws.Range("E57").Value = 49500
$ws.Range("D58").Value = 700
$ws.Range("E58").Value = 400
$ws.Range("D59").Value = 61500
$ws.Range("E59").Value = 44600
$ws.Range("D60").Value = 108300
$ws.Range("E60").Value = 94500
$ws.Range("D61").Value = 427300
$ws.Range("E61").Value = 115300
$ws.Range("D62").Value = 10800
$ws.Range("E62").Value = 6000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 546300
$ws.Range("E66").Value = 215800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -147100
$ws.Range("E72").Value = -69700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 594800
$ws.Range("E76").Value = 490600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -77300
$ws.Range("E81").Value = 13700
$ws.Range("D83").Value = 18200
$ws.Range("E83").Value = 15500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 38800
$ws.Range("E89").Value = 25600
$ws.Range("D91").Value = -17100
$ws.Range("E91").Value = -11500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -363800
$ws.Range("E94").Value = -11000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 302100
$ws.Range("E100").Value = 1100
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -22900
$ws.Range("E102").Value = 15700
